$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
# Row 33
$ws.Cells.Item(33, 8).Value = 158  # H33: was 158.69565
$ws.Cells.Item(33, 9).Value = 131.70589  # I33: was 139.8421
$ws.Cells.Item(33, 10).Value = 307  # J33: was 248.25
$ws.Cells.Item(33, 11).Value = 131.70589  # K33: was 139.8421
$ws.Cells.Item(33, 12).Value = 307  # L33: was 248.25
$ws.Cells.Item(33, 13).Value = 97.29410999999999  # M33: was 89.15790000000001
$ws.Cells.Item(33, 14).Value = -765  # N33: was -706.25

# Row 40
$ws.Cells.Item(40, 8).Value = 3524.5  # H40: was 3459.6
$ws.Cells.Item(40, 10).Value = 3632.6667  # J40: was 3524.5
$ws.Cells.Item(40, 12).Value = 3632.6667  # L40: was 3524.5
$ws.Cells.Item(40, 14).Value = -3982.6667  # N40: was -3874.5

# Row 62
$ws.Cells.Item(62, 8).Value = 3470.5  # H62: was 3470.75
$ws.Cells.Item(62, 9).Value = 2998  # I62: was 2998.3333
$ws.Cells.Item(62, 11).Value = 2998  # K62: was 2998.3333
$ws.Cells.Item(62, 13).Value = -2374  # M62: was -2374.3333

# Row 65
$ws.Cells.Item(65, 8).Value = 3470.5  # H65: was 3470.75
$ws.Cells.Item(65, 9).Value = 2998  # I65: was 2998.3333
$ws.Cells.Item(65, 11).Value = 14990  # K65: was 14991.6665
$ws.Cells.Item(65, 13).Value = -11870  # M65: was -11871.6665

# Row 76
$ws.Cells.Item(76, 8).Value = 4042.2856  # H76: was 4099.7144
$ws.Cells.Item(76, 9).Value = 4179.4  # I76: was 4259.8
$ws.Cells.Item(76, 11).Value = 4179.4  # K76: was 4259.8
$ws.Cells.Item(76, 13).Value = -3864.4  # M76: was -3944.8

# Row 79
$ws.Cells.Item(79, 8).Value = 4042.2856  # H79: was 4099.7144
$ws.Cells.Item(79, 9).Value = 4179.4  # I79: was 4259.8
$ws.Cells.Item(79, 11).Value = 4179.4  # K79: was 4259.8
$ws.Cells.Item(79, 13).Value = -3087.4  # M79: was -3167.8

# Row 111
$ws.Cells.Item(111, 8).Value = 2482.5  # H111: was 2650
$ws.Cells.Item(111, 10).Value = 1993.3334  # J111: was 2000
$ws.Cells.Item(111, 12).Value = 5980.0002  # L111: was 6000
$ws.Cells.Item(111, 14).Value = -12114.0002  # N111: was -12134

# Row 113
$ws.Cells.Item(113, 8).Value = 19287  # H113: was 25001
$ws.Cells.Item(113, 9).Value = 23000.8  # I113: was 35000
$ws.Cells.Item(113, 11).Value = 23000.8  # K113: was 35000
$ws.Cells.Item(113, 13).Value = -19746.8  # M113: was -31746

# Row 127
$ws.Cells.Item(127, 8).Value = 459.8  # H127: was 449
$ws.Cells.Item(127, 9).Value = 499.66666  # I127: was 449
$ws.Cells.Item(127, 10).Value = 400  # J127: was 0
$ws.Cells.Item(127, 11).Value = 1498.99998  # K127: was 1347
$ws.Cells.Item(127, 12).Value = 1200  # L127: was 0
$ws.Cells.Item(127, 13).Value = 3461.00002  # M127: was 3613
$ws.Cells.Item(127, 14).Value = -11120  # N127: was None

# Row 137
$ws.Cells.Item(137, 8).Value = 1331.5652  # H137: was 1334.1666
$ws.Cells.Item(137, 9).Value = 1356.5  # I137: was 1321.5385
$ws.Cells.Item(137, 10).Value = 1304.3636  # J137: was 1349.091
$ws.Cells.Item(137, 11).Value = 4069.5  # K137: was 3964.6155
$ws.Cells.Item(137, 12).Value = 3913.0908  # L137: was 4047.273
$ws.Cells.Item(137, 13).Value = -1519.5  # M137: was -1414.6155
$ws.Cells.Item(137, 14).Value = -9013.0908  # N137: was -9147.272999999999

# --- Sheet: ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
# Row 2
$ws.Cells.Item(2, 8).Value = 1614.2778  # H2: was 1555.5264
$ws.Cells.Item(2, 9).Value = 909.5  # I2: was 872.0909
$ws.Cells.Item(2, 11).Value = 909.5  # K2: was 872.0909
$ws.Cells.Item(2, 13).Value = -796.5  # M2: was -759.0909

# Row 32
$ws.Cells.Item(32, 8).Value = 3712.875  # H32: was 3486.2903
$ws.Cells.Item(32, 9).Value = 2647.6296  # I32: was 2336.5
$ws.Cells.Item(32, 11).Value = 2647.6296  # K32: was 2336.5
$ws.Cells.Item(32, 13).Value = -2360.6296  # M32: was -2049.5

# Row 45
$ws.Cells.Item(45, 8).Value = 3416.9048  # H45: was 3350.182
$ws.Cells.Item(45, 10).Value = 3922.25  # J45: was 3806.1765
$ws.Cells.Item(45, 12).Value = 3922.25  # L45: was 3806.1765
$ws.Cells.Item(45, 14).Value = -4676.25  # N45: was -4560.1765

# Row 74
$ws.Cells.Item(74, 8).Value = 804.2222  # H74: was 802.6
$ws.Cells.Item(74, 9).Value = 804.2222  # I74: was 802.6
$ws.Cells.Item(74, 11).Value = 804.2222  # K74: was 802.6
$ws.Cells.Item(74, 13).Value = 69.77779999999996  # M74: was 71.39999999999998

# Row 77
$ws.Cells.Item(77, 8).Value = 804.2222  # H77: was 802.6
$ws.Cells.Item(77, 9).Value = 804.2222  # I77: was 802.6
$ws.Cells.Item(77, 11).Value = 4021.111  # K77: was 4013
$ws.Cells.Item(77, 13).Value = 346.8889999999997  # M77: was 355

# Row 110
$ws.Cells.Item(110, 8).Value = 934.8889  # H110: was 964.375
$ws.Cells.Item(110, 9).Value = 839.6667  # I110: was 867.8
$ws.Cells.Item(110, 11).Value = 839.6667  # K110: was 867.8
$ws.Cells.Item(110, 13).Value = 1205.3333  # M110: was 1177.2

# Row 116
$ws.Cells.Item(116, 8).Value = 1614.2778  # H116: was 1555.5264
$ws.Cells.Item(116, 9).Value = 909.5  # I116: was 872.0909
$ws.Cells.Item(116, 11).Value = 909.5  # K116: was 872.0909
$ws.Cells.Item(116, 13).Value = 1384.5  # M116: was 1421.9091

# Row 130
$ws.Cells.Item(130, 8).Value = 9997  # H130: was 9999.333000000001
$ws.Cells.Item(130, 10).Value = 9997  # J130: was 9999.333000000001
$ws.Cells.Item(130, 12).Value = 9997  # L130: was 9999.333000000001
$ws.Cells.Item(130, 14).Value = -20037  # N130: was -20039.333

# --- Sheet: BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)
# Row 3
$ws.Cells.Item(3, 8).Value = 1614.2778  # H3: was 1555.5264
$ws.Cells.Item(3, 9).Value = 909.5  # I3: was 872.0909
$ws.Cells.Item(3, 11).Value = 909.5  # K3: was 872.0909
$ws.Cells.Item(3, 13).Value = -795.5  # M3: was -758.0909

# Row 94
$ws.Cells.Item(94, 8).Value = 2488.5  # H94: was 3282.5715
$ws.Cells.Item(94, 9).Value = 982  # I94: was 999.5
$ws.Cells.Item(94, 10).Value = 3995  # J94: was 4195.8
$ws.Cells.Item(94, 11).Value = 982  # K94: was 999.5
$ws.Cells.Item(94, 12).Value = 3995  # L94: was 4195.8
$ws.Cells.Item(94, 13).Value = -531  # M94: was -548.5
$ws.Cells.Item(94, 14).Value = -4897  # N94: was -5097.8

# Row 107
$ws.Cells.Item(107, 8).Value = 999.5  # H107: was 833
$ws.Cells.Item(107, 9).Value = 999.5  # I107: was 833
$ws.Cells.Item(107, 11).Value = 999.5  # K107: was 833
$ws.Cells.Item(107, 13).Value = 920.5  # M107: was 1087

# Row 135
$ws.Cells.Item(135, 8).Value = 39999.5  # H135: was 40000
$ws.Cells.Item(135, 10).Value = 39999.332  # J135: was 40000
$ws.Cells.Item(135, 12).Value = 39999.332  # L135: was 40000
$ws.Cells.Item(135, 14).Value = -50139.332  # N135: was -50140

# --- Sheet: CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
# Row 22
$ws.Cells.Item(22, 8).Value = 1512.5  # H22: was 1841.6666
$ws.Cells.Item(22, 9).Value = 1512.5  # I22: was 1683.3334
$ws.Cells.Item(22, 10).Value = 0  # J22: was 2000
$ws.Cells.Item(22, 11).Value = 1512.5  # K22: was 1683.3334
$ws.Cells.Item(22, 12).Value = 0  # L22: was 2000
$ws.Cells.Item(22, 13).Value = -1162.5  # M22: was -1333.3334
$ws.Cells.Item(22, 14).ClearContents()  # N22: was -2700, removed

# Row 31
$ws.Cells.Item(31, 8).Value = 2117.7144  # H31: was 1744
$ws.Cells.Item(31, 9).Value = 1441.6666  # I31: was 1247.2
$ws.Cells.Item(31, 10).Value = 2624.75  # J31: was 2158
$ws.Cells.Item(31, 11).Value = 1441.6666  # K31: was 1247.2
$ws.Cells.Item(31, 12).Value = 2624.75  # L31: was 2158
$ws.Cells.Item(31, 13).Value = -1146.6666  # M31: was -952.2
$ws.Cells.Item(31, 14).Value = -3214.75  # N31: was -2748

# Row 34
$ws.Cells.Item(34, 8).Value = 2117.7144  # H34: was 1744
$ws.Cells.Item(34, 9).Value = 1441.6666  # I34: was 1247.2
$ws.Cells.Item(34, 10).Value = 2624.75  # J34: was 2158
$ws.Cells.Item(34, 11).Value = 1441.6666  # K34: was 1247.2
$ws.Cells.Item(34, 12).Value = 2624.75  # L34: was 2158
$ws.Cells.Item(34, 13).Value = -1239.6666  # M34: was -1045.2
$ws.Cells.Item(34, 14).Value = -3028.75  # N34: was -2562

# Row 99
$ws.Cells.Item(99, 8).Value = 6374.6665  # H99: was 6789.9
$ws.Cells.Item(99, 9).Value = 6317.8184  # I99: was 6766.5557
$ws.Cells.Item(99, 11).Value = 6317.8184  # K99: was 6766.5557
$ws.Cells.Item(99, 13).Value = -4819.8184  # M99: was -5268.5557

# Row 122
$ws.Cells.Item(122, 8).Value = 2832.25  # H122: was 2784
$ws.Cells.Item(122, 9).Value = 2864  # I122: was 2805.75
$ws.Cells.Item(122, 11).Value = 8592  # K122: was 8417.25
$ws.Cells.Item(122, 13).Value = -6142  # M122: was -5967.25

# Row 126
$ws.Cells.Item(126, 8).Value = 6374.6665  # H126: was 6789.9
$ws.Cells.Item(126, 9).Value = 6317.8184  # I126: was 6766.5557
$ws.Cells.Item(126, 11).Value = 18953.4552  # K126: was 20299.6671
$ws.Cells.Item(126, 13).Value = -16483.4552  # M126: was -17829.6671

# Row 132
$ws.Cells.Item(132, 8).Value = 2770.7144  # H132: was 2437.6667
$ws.Cells.Item(132, 9).Value = 2770.7144  # I132: was 2437.6667
$ws.Cells.Item(132, 11).Value = 8312.143199999999  # K132: was 7313.000100000001
$ws.Cells.Item(132, 13).Value = -5782.143199999999  # M132: was -4783.000100000001

# --- Sheet: CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)
# Row 39
$ws.Cells.Item(39, 8).Value = 2000  # H39: was 3000
$ws.Cells.Item(39, 10).Value = 2000  # J39: was 3000
$ws.Cells.Item(39, 12).Value = 6000  # L39: was 9000
$ws.Cells.Item(39, 14).Value = -6588  # N39: was -9588

# Row 55
$ws.Cells.Item(55, 8).Value = 10004  # H55: was 4145.143
$ws.Cells.Item(55, 9).Value = 10004  # I55: was 2504
$ws.Cells.Item(55, 10).Value = 0  # J55: was 6333.3335
$ws.Cells.Item(55, 11).Value = 30012  # K55: was 7512
$ws.Cells.Item(55, 12).Value = 0  # L55: was 19000.0005
$ws.Cells.Item(55, 13).Value = -29835  # M55: was -7335
$ws.Cells.Item(55, 14).ClearContents()  # N55: was -19354.0005, removed

# Row 125
$ws.Cells.Item(125, 8).Value = 10500  # H125: was 25000
$ws.Cells.Item(125, 9).Value = 1000  # I125: was 0
$ws.Cells.Item(125, 10).Value = 20000  # J125: was 25000
$ws.Cells.Item(125, 11).Value = 3000  # K125: was 0
$ws.Cells.Item(125, 12).Value = 60000  # L125: was 75000
$ws.Cells.Item(125, 13).Value = 1920  # M125: was None
$ws.Cells.Item(125, 14).Value = -69840  # N125: was -84840

# --- Sheet: GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
# Row 3
$ws.Cells.Item(3, 8).Value = 3000  # H3: was 5000000
$ws.Cells.Item(3, 9).Value = 3000  # I3: was 5000000
$ws.Cells.Item(3, 11).Value = 3000  # K3: was 5000000
$ws.Cells.Item(3, 13).Value = -2884  # M3: was -4999884

# Row 12
$ws.Cells.Item(12, 8).Value = 28000000  # H12: was 0
$ws.Cells.Item(12, 9).Value = 28000000  # I12: was 0
$ws.Cells.Item(12, 11).Value = 28000000  # K12: was 0
$ws.Cells.Item(12, 13).Value = -27999860  # M12: was None

# Row 122
$ws.Cells.Item(122, 8).Value = 2170.9473  # H122: was 2055.6667
$ws.Cells.Item(122, 9).Value = 2014.1177  # I122: was 1903.2106
$ws.Cells.Item(122, 11).Value = 6042.3531  # K122: was 5709.6318
$ws.Cells.Item(122, 13).Value = -3592.3531  # M122: was -3259.6318

# Row 132
$ws.Cells.Item(132, 8).Value = 2294.4285  # H132: was 2618.182
$ws.Cells.Item(132, 9).Value = 2355.5386  # I132: was 2618.182
$ws.Cells.Item(132, 10).Value = 1500  # J132: was 0
$ws.Cells.Item(132, 11).Value = 7066.6158  # K132: was 7854.545999999999
$ws.Cells.Item(132, 12).Value = 4500  # L132: was 0
$ws.Cells.Item(132, 13).Value = -4536.6158  # M132: was -5324.545999999999
$ws.Cells.Item(132, 14).Value = -9560  # N132: was None

# --- Sheet: LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
# Row 29
$ws.Cells.Item(29, 8).Value = 0  # H29: was 50000000
$ws.Cells.Item(29, 9).Value = 0  # I29: was 50000000
$ws.Cells.Item(29, 11).Value = 0  # K29: was 50000000
$ws.Cells.Item(29, 13).ClearContents()  # M29: was -49999705, removed

# Row 61
$ws.Cells.Item(61, 8).Value = 1457.9166  # H61: was 1771.4286
$ws.Cells.Item(61, 9).Value = 1223.125  # I61: was 1497.5
$ws.Cells.Item(61, 10).Value = 1927.5  # J61: was 2136.6667
$ws.Cells.Item(61, 11).Value = 1223.125  # K61: was 1497.5
$ws.Cells.Item(61, 12).Value = 1927.5  # L61: was 2136.6667
$ws.Cells.Item(61, 13).Value = -1021.125  # M61: was -1295.5
$ws.Cells.Item(61, 14).Value = -2331.5  # N61: was -2540.6667

# Row 93
$ws.Cells.Item(93, 8).Value = 964.8333  # H93: was 867.1818
$ws.Cells.Item(93, 9).Value = 917.8  # I93: was 827.1429000000001
$ws.Cells.Item(93, 10).Value = 1200  # J93: was 937.25
$ws.Cells.Item(93, 11).Value = 917.8  # K93: was 827.1429000000001
$ws.Cells.Item(93, 12).Value = 1200  # L93: was 937.25
$ws.Cells.Item(93, 13).Value = 330.2  # M93: was 420.8570999999999
$ws.Cells.Item(93, 14).Value = -3696  # N93: was -3433.25

# Row 113
$ws.Cells.Item(113, 8).Value = 1457.9166  # H113: was 1771.4286
$ws.Cells.Item(113, 9).Value = 1223.125  # I113: was 1497.5
$ws.Cells.Item(113, 10).Value = 1927.5  # J113: was 2136.6667
$ws.Cells.Item(113, 11).Value = 1223.125  # K113: was 1497.5
$ws.Cells.Item(113, 12).Value = 1927.5  # L113: was 2136.6667
$ws.Cells.Item(113, 13).Value = 946.875  # M113: was 672.5
$ws.Cells.Item(113, 14).Value = -6267.5  # N113: was -6476.6667

# Row 122
$ws.Cells.Item(122, 8).Value = 1567.3334  # H122: was 1411.4445
$ws.Cells.Item(122, 9).Value = 1567.3334  # I122: was 1411.4445
$ws.Cells.Item(122, 11).Value = 4702.0002  # K122: was 4234.333500000001
$ws.Cells.Item(122, 13).Value = -2252.0002  # M122: was -1784.333500000001

# Row 136
$ws.Cells.Item(136, 8).Value = 2271.85  # H136: was 2439.5264
$ws.Cells.Item(136, 9).Value = 2331  # I136: was 2390.5334
$ws.Cells.Item(136, 10).Value = 2094.4  # J136: was 2623.25
$ws.Cells.Item(136, 11).Value = 6993  # K136: was 7171.600199999999
$ws.Cells.Item(136, 12).Value = 6283.200000000001  # L136: was 7869.75
$ws.Cells.Item(136, 13).Value = -4443  # M136: was -4621.600199999999
$ws.Cells.Item(136, 14).Value = -11383.2  # N136: was -12969.75

# --- Sheet: WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)
# Row 81
$ws.Cells.Item(81, 8).Value = 5489.8887  # H81: was 5492.5557
$ws.Cells.Item(81, 9).Value = 5629.857  # I81: was 5738.8335
$ws.Cells.Item(81, 11).Value = 11259.714  # K81: was 11477.667
$ws.Cells.Item(81, 13).Value = -10198.714  # M81: was -10416.667

# Row 84
$ws.Cells.Item(84, 8).Value = 5489.8887  # H84: was 5492.5557
$ws.Cells.Item(84, 9).Value = 5629.857  # I84: was 5738.8335
$ws.Cells.Item(84, 11).Value = 56298.57  # K84: was 57388.335
$ws.Cells.Item(84, 13).Value = -50994.57  # M84: was -52084.335

# Row 96
$ws.Cells.Item(96, 8).Value = 1829.9333  # H96: was 2050
$ws.Cells.Item(96, 9).Value = 1879.9  # I96: was 1900
$ws.Cells.Item(96, 10).Value = 1730  # J96: was 2275
$ws.Cells.Item(96, 11).Value = 1879.9  # K96: was 1900
$ws.Cells.Item(96, 12).Value = 1730  # L96: was 2275
$ws.Cells.Item(96, 13).Value = -506.9000000000001  # M96: was -527
$ws.Cells.Item(96, 14).Value = -4476  # N96: was -5021

# Row 110
$ws.Cells.Item(110, 8).Value = 0  # H110: was 80644
$ws.Cells.Item(110, 10).Value = 0  # J110: was 80644
$ws.Cells.Item(110, 12).Value = 0  # L110: was 80644
$ws.Cells.Item(110, 14).ClearContents()  # N110: was -88824, removed

# Row 122
$ws.Cells.Item(122, 8).Value = 3112.25  # H122: was 1981.875
$ws.Cells.Item(122, 9).Value = 2725  # I122: was 1476
$ws.Cells.Item(122, 11).Value = 8175  # K122: was 4428
$ws.Cells.Item(122, 13).Value = -5725  # M122: was -1978

# Row 126
$ws.Cells.Item(126, 8).Value = 2608.1428  # H126: was 1878.3636
$ws.Cells.Item(126, 9).Value = 2608.1428  # I126: was 1878.3636
$ws.Cells.Item(126, 11).Value = 7824.428400000001  # K126: was 5635.0908
$ws.Cells.Item(126, 13).Value = -5354.428400000001  # M126: was -3165.0908

# Row 132
$ws.Cells.Item(132, 8).Value = 2225.6667  # H132: was 2182
$ws.Cells.Item(132, 9).Value = 1920.8  # I132: was 1920.6666
$ws.Cells.Item(132, 11).Value = 5762.4  # K132: was 5761.9998
$ws.Cells.Item(132, 13).Value = -3232.4  # M132: was -3231.9998

# Row 136
$ws.Cells.Item(136, 8).Value = 686.1923  # H136: was 674.96295
$ws.Cells.Item(136, 9).Value = 738.6957  # I136: was 723.875
$ws.Cells.Item(136, 11).Value = 2216.0871  # K136: was 2171.625
$ws.Cells.Item(136, 13).Value = 333.9129000000003  # M136: was 378.375

